$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")

# The row describing "my_last_frame" is row 10 (A10:F10). Remove it
# (clear contents) everywhere, matching "remove my_last_frame from everywhere".
$rng = $ws.Range("A10:F10")
$rng.Clear()
$rng.Select()
